$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 ("citizen") test-data row got its credentials refreshed as part of
# the grievances/council-management refactor: old phone-number id
# "9036544535" / password "akhi2506" are replaced with a new id
# "7259960764" using the already-shared "kurnool_eGov@123" password.
$ws.Range("B26").Value = "7259960764"
$ws.Range("C26").Value = "kurnool_eGov@123"

# Move the active selection to C33 (the sheet was scrolled further down
# while editing this data).
[void]$ws.Range("C33").Select()
